$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:B1").Copy()
$ws.Range("A3:B3").PasteSpecial(-4122)

$ws.Range("A3").Value = "passive income"
$ws.Range("B3").Value = "passive.income.nadi.myfirstdrawermenuproject"

$ws.Rows.Item(3).RowHeight = 24

$ws.Range("A3").Select()
